$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in C4 (Points Completed for Sprint 3)
$ws.Range("C4").Value = 38

# Add new row for Sprint 4
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 25

# Update selection to match the new active cell
$ws.Range("B5").Select()
